$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("link_buget_example 1")
$ws.Activate()

# ---------------------------------------------------------------------------
# 1. Link budget inputs: E3/E9/E12/E15 become plain literal values (previously
#    E9/E12/E15 were small addition formulas). Everything downstream
#    (A3,A9,A12,A15,A25,B25,A30,B30,C30,B35) recalculates automatically.
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 108.99
$ws.Range("E9").Value = 110.86
$ws.Range("E12").Value = 110.86
$ws.Range("E15").Value = 120.98

# SE (bps/Hz) input changes from 5.5547 to 1.4766
$ws.Range("A35").Value = 1.4766

# ---------------------------------------------------------------------------
# 2. Rebuild the "Capacity" block (old rows 38-42) into two stacked blocks:
#    rows 37-42 (Capacity (bps) / Capacity (Mbps)) and rows 44-49
#    (Capacity (bps/km2) / Capacity (Mbps/km2)). Clear the old block first.
# ---------------------------------------------------------------------------
$ws.Range("A38:H42").Clear()

# -- Row 37: "Capacity" section header (same style as row 33/18/etc.) --------
$ws.Range("A37").Value = "Capacity"
$ws.Range("A37").Font.Bold = $true
$ws.Range("A37").Font.Size = 12
$ws.Rows.Item(37).RowHeight = 16.5

# -- Row 38: column headers ---------------------------------------------------
$ws.Range("A38").Value = "Capacity (bps)"
$ws.Range("B38").Value = "SE"
$ws.Range("C38").Value = "Bandwidth (Hz)"
$ws.Range("D38").Value = "Area (km^2)"
$ws.Range("A38:D38").Font.Bold = $true
$ws.Range("A38:D38").HorizontalAlignment = -4108

# -- Row 39: capacity (bps) formula row --------------------------------------
$ws.Range("A39").Formula = "=B39*C39"
$ws.Range("B39").Formula = "=A35"
$ws.Range("C39").Formula = "=10000000"
$ws.Range("D39").Value = 0.2165063
$ws.Rows.Item(39).RowHeight = 15.75

# -- Row 40: blank thick-bottom spacer row -----------------------------------
$ws.Rows.Item(40).RowHeight = 15.75

# -- Row 41: "Capacity (Mbps)" label with boxed top border -------------------
$ws.Range("A41").Value = "Capacity (Mbps)"
$ws.Range("A41").Font.Bold = $true
$ws.Range("A41").HorizontalAlignment = -4108

# -- Row 42: capacity (Mbps) value with boxed bottom border ------------------
$ws.Range("A42").Formula = "=A39/1000000"
$ws.Range("A42").NumberFormat = "_-* #,##0.00_-;\-* #,##0.00_-;_-* " + [char]34 + "-" + [char]34 + "??_-;_-@_-"
$ws.Range("A42").Font.Color = $ws.Range("D25").Font.Color
$ws.Rows.Item(42).RowHeight = 15.75

# Box borders around A41:A42 (medium outline, thin line between label/value)
$ws.Range("A41:A42").Borders.Item(7).LineStyle = 1
$ws.Range("A41:A42").Borders.Item(7).Weight = -4138
$ws.Range("A41:A42").Borders.Item(10).LineStyle = 1
$ws.Range("A41:A42").Borders.Item(10).Weight = -4138
$ws.Range("A41").Borders.Item(8).LineStyle = 1
$ws.Range("A41").Borders.Item(8).Weight = -4138
$ws.Range("A41").Borders.Item(9).LineStyle = 1
$ws.Range("A41").Borders.Item(9).Weight = 2
$ws.Range("A42").Borders.Item(9).LineStyle = 1
$ws.Range("A42").Borders.Item(9).Weight = -4138

# -- Row 44: "Capacity" section header (second block) ------------------------
$ws.Range("A44").Value = "Capacity"
$ws.Range("A44").Font.Bold = $true
$ws.Range("A44").Font.Size = 12
$ws.Rows.Item(44).RowHeight = 16.5

# -- Row 45: column headers ---------------------------------------------------
$ws.Range("A45").Value = "Capacity (bps/km2)"
$ws.Range("B45").Value = "SE"
$ws.Range("C45").Value = "Bandwidth (Hz)"
$ws.Range("D45").Value = "Area (km^2)"
$ws.Range("A45:D45").Font.Bold = $true
$ws.Range("A45:D45").HorizontalAlignment = -4108

# -- Row 46: capacity (bps/km2) formula row ----------------------------------
$ws.Range("A46").Formula = "=B46*C46/D46"
$ws.Range("B46").Formula = "=A35"
$ws.Range("C46").Formula = "=E20"
$ws.Range("D46").Value = 0.2165063
$ws.Rows.Item(46).RowHeight = 15.75

# -- Row 47: blank thick-bottom spacer row -----------------------------------
$ws.Rows.Item(47).RowHeight = 15.75

# -- Row 48: "Capacity (Mbps/km2)" label with boxed top border ---------------
$ws.Range("A48").Value = "Capacity (Mbps/km2)"
$ws.Range("A48").Font.Bold = $true
$ws.Range("A48").HorizontalAlignment = -4108

# -- Row 49: capacity (Mbps/km2) value with boxed bottom border --------------
$ws.Range("A49").Formula = "=A46/1000000"
$ws.Range("A49").NumberFormat = "_-* #,##0.00_-;\-* #,##0.00_-;_-* " + [char]34 + "-" + [char]34 + "??_-;_-@_-"
$ws.Range("A49").Font.Color = $ws.Range("D25").Font.Color
$ws.Rows.Item(49).RowHeight = 15.75

$ws.Range("A48:A49").Borders.Item(7).LineStyle = 1
$ws.Range("A48:A49").Borders.Item(7).Weight = -4138
$ws.Range("A48:A49").Borders.Item(10).LineStyle = 1
$ws.Range("A48:A49").Borders.Item(10).Weight = -4138
$ws.Range("A48").Borders.Item(8).LineStyle = 1
$ws.Range("A48").Borders.Item(8).Weight = -4138
$ws.Range("A48").Borders.Item(9).LineStyle = 1
$ws.Range("A48").Borders.Item(9).Weight = 2
$ws.Range("A49").Borders.Item(9).LineStyle = 1
$ws.Range("A49").Borders.Item(9).Weight = -4138

# ---------------------------------------------------------------------------
# 3. View state: zoom + scroll position + active selection.
# ---------------------------------------------------------------------------
$ws.Range("A8").Select()
$excel.ActiveWindow.ScrollRow = 8
$excel.ActiveWindow.ScrollColumn = 1
$excel.ActiveWindow.Zoom = 80
$ws.Range("B51").Select()
